$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the submitter email value (shared string change)
$ws.Range("AK2").Value = "sindhube19.data@gmail.com"
